$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at position 932, pushing the existing data (and all rows below)
# down by two rows. This mirrors the diff: old rows 932-995 become new rows 934-997,
# and two brand-new rows of weekly data are inserted at 932-933.
$ws.Rows.Item(932).Insert()
$ws.Rows.Item(932).Insert()

# Row 932: new weekly entry (Primera, $/caja 36 atados)
$ws.Cells.Item(932, 1).Value = 6
$ws.Cells.Item(932, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(932, 3).Value = "Metropolitana"
$ws.Cells.Item(932, 4).Value = 44826
$ws.Cells.Item(932, 5).Value = 13
$ws.Cells.Item(932, 6).Value = 100112040
$ws.Cells.Item(932, 7).Value = "Cilantro"
$ws.Cells.Item(932, 8).Value = "Sin especificar"
$ws.Cells.Item(932, 9).Value = "Primera"
$ws.Cells.Item(932, 10).Value = 630
$ws.Cells.Item(932, 11).Value = 4500
$ws.Cells.Item(932, 12).Value = 5000
$ws.Cells.Item(932, 13).Value = 4714
$ws.Cells.Item(932, 14).Value = "`$/caja 36 atados"
$ws.Cells.Item(932, 15).Value = "Región Metropolitana"
$ws.Cells.Item(932, 16).Value = 131
$ws.Cells.Item(932, 17).Value = 36
$ws.Cells.Item(932, 18).Value = "Hortaliza"

# Row 933: new weekly entry (Primera, $/docena de atados)
$ws.Cells.Item(933, 1).Value = 6
$ws.Cells.Item(933, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(933, 3).Value = "Metropolitana"
$ws.Cells.Item(933, 4).Value = 44826
$ws.Cells.Item(933, 5).Value = 13
$ws.Cells.Item(933, 6).Value = 100112040
$ws.Cells.Item(933, 7).Value = "Cilantro"
$ws.Cells.Item(933, 8).Value = "Sin especificar"
$ws.Cells.Item(933, 9).Value = "Primera"
$ws.Cells.Item(933, 10).Value = 450
$ws.Cells.Item(933, 11).Value = 7500
$ws.Cells.Item(933, 12).Value = 8000
$ws.Cells.Item(933, 13).Value = 7711
$ws.Cells.Item(933, 14).Value = "`$/docena de atados"
$ws.Cells.Item(933, 15).Value = "Región Metropolitana"
$ws.Cells.Item(933, 16).Value = 2570
$ws.Cells.Item(933, 17).Value = 3
$ws.Cells.Item(933, 18).Value = "Hortaliza"
